$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Remove the "GDP growth" row (old row 4) - folded away, data no longer
#    tracked separately from "GDP per capita".
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).Delete()

# ---------------------------------------------------------------------------
# 2) "Agricultural sector value of GDP" row (now row 5): add a Details entry.
# ---------------------------------------------------------------------------
$ws.Range("E5").Value = "Proportion of national GDP"

# ---------------------------------------------------------------------------
# 3) "development flows to environment" row (now row 7): fix capitalisation.
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "Development flows to environment"

# ---------------------------------------------------------------------------
# 4) Remove the "Agricultural Raw Materials" / IMF commodity-price row
#    (now row 10) - replaced by the more detailed commodity price rows below.
# ---------------------------------------------------------------------------
$ws.Rows.Item(10).Delete()

# ---------------------------------------------------------------------------
# 5) Append the new commodity-price rows (rice / corn / rubber / sugar).
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = "Price of rice"
$ws.Range("B13").Value = "USD/ton"
$ws.Range("C13").Value = "Global"
$ws.Range("D13").Value = "World Bank"
$ws.Range("E13").Value = "Median annual global market price of rice"
$ws.Range("E13").WrapText = $true

$ws.Range("A14").Value = "Price of corn"
$ws.Range("B14").Value = "USD/ton"
$ws.Range("C14").Value = "Global"
$ws.Range("D14").Value = "World Bank"
$ws.Range("E14").Value = "Annual global market price of corn"
$ws.Range("E14").WrapText = $true

$ws.Range("A15").Value = "Price of rubber"
$ws.Range("B15").Value = "USD/ton"
$ws.Range("C15").Value = "Regional"
$ws.Range("D15").Value = "RASCE"
$ws.Range("E15").Value = "Monthly regional market value of rubber on the Singapore Exchange"
$ws.Range("E15").WrapText = $true

$ws.Range("A16").Value = "Price of sugar"
$ws.Range("B16").Value = "USD/ton"
$ws.Range("C16").Value = "Global"
$ws.Range("D16").Value = "World Bank"
$ws.Range("E16").Value = "Annual global market price of sugar"
$ws.Range("E16").WrapText = $true

# ---------------------------------------------------------------------------
# 6) New "Producer prices" section with four follow-up rows (left for later
#    detail - source/units/etc. still to be filled in by the author).
# ---------------------------------------------------------------------------
$ws.Range("A17").Value = "Producer prices"
$ws.Range("A17").Font.Italic = $true

$ws.Range("A18").Value = "Producer price of rubber"
$ws.Range("A19").Value = "Producer price of cassava"
$ws.Range("A20").Value = "Producer price of corn"
$ws.Range("A21").Value = "Producer price of sugar"

# ---------------------------------------------------------------------------
# 7) Widen columns D (Source) and E (Details) to accommodate the new content.
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 20.333333333333332
$ws.Columns.Item(5).ColumnWidth = 60.166666666666664
